$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (columns B:S) down by one row, from the bottom up,
# to make room for the new company row that is inserted at row 2.
# Column A (the numeric index 0,1,2,...) is intentionally left untouched --
# each row keeps the A value that already belonged to that row position.
for ($r = 30; $r -ge 2; $r--) {
    $destRow = $r + 1
    $srcRange = $ws.Range("B" + $r + ":S" + $r)
    $destRange = $ws.Range("B" + $destRow + ":S" + $destRow)
    $destRange.Value2 = $srcRange.Value2
}

# Give the new last row (31) the same numbering-column format as the rest of column A,
# then set its value to continue the 0-based sequence (29).
$ws.Range("A30").Copy($ws.Range("A31"))
$ws.Range("A31").Value2 = 29

# Populate the brand-new row 2 with the new company data.
$ws.Range("B2").Value2 = '新康众有限公司（江苏康众汽配）'
$ws.Range("C2").Value2 = '江苏省南京市建邺区'
$ws.Range("D2").Value2 = '技术中心-供应链'
$ws.Range("E2").Value2 = 'Java'
$ws.Range("F2").Value2 = '9:00-18：00'
$ws.Range("G2").Value2 = '1.5h'
$ws.Range("H2").Value2 = '不强制加班，但是你没法早走，经常晚上开会'
$ws.Range("I2").Value2 = '个人1500+公司1500'
$ws.Range("J2").Value2 = 'hr说是3个月，实则0'
$ws.Range("K2").Value2 = '试用期3个月；不打折。'
$ws.Range("L2").Value2 = '网吧工位，电脑自己买，公司最多报销1w，分36个月返给你'
$ws.Range("M2").Value2 = '入职满3年有三天，且没有陪产假这些。'
$ws.Range("N2").Value2 = '钉钉打卡，领导会看监控以防你打卡了但人没到'
$ws.Range("O2").Value2 = '公司没有任何福利，老板很抠门。从大领导到小领导都只会pua。直系领导每天啥都不干，只盯着你不让你闲着，要求24小时保持工作状态'
$ws.Range("P2").Value2 = ''
$ws.Range("Q2").Value2 = '2022-06-23 10:00:12'
$ws.Range("R2").Value2 = ''
$ws.Range("S2").Value2 = ''

Write-Host "Row insert complete"
